# Rephrasing suggestion of EmptyCatch check.
#
# Checklist.xlsx has two sheets: "Workflow" (checklist of workflow-level
# checks) and "Project". On the "Workflow" sheet, row 5 is the
# "Empty Catch block" check, and column G holds the "Suggestion" text.
# Only the wording of that suggestion changes.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Workflow")

$ws.Range("G5").Value = "If no action is taken to handle the exception, consider including at least a Log Message activity and Rethrow it."

# Leave the cursor where the author's edit left it.
$ws.Activate()
$ws.Range("G6").Select()
